# Update column F (dSF) values for specific rows to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 6
    10 = 1
    11 = 0
    14 = 2
    17 = 3
    23 = -2
    25 = -4
    28 = 0
    30 = 1
    35 = 1
    37 = 5
    38 = -5
    41 = 4
    47 = 2
    51 = -2
    54 = -3
    58 = 1
    61 = -5
    73 = -3
    77 = -2
    79 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
